# CAPI-Licences.xlsx: "Added license for CSS."
#
# 1) The ClearCanvas row's Licence cell was corrected from "GNU" to "GPLv3"
#    (this also makes the old "GNU" shared string unused/dropped on save).
# 2) A new row (11) was appended documenting the "Now UI Dashboard" CSS
#    template used by WebService/assets, with its licence, link and a
#    comment.
# 3) Hyperlinks were attached to the Link column for the new row (D11) and
#    for the ClearCanvas row (D2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the ClearCanvas licence -----------------------------------
$ws.Range("C2").Value = "GPLv3"

# --- 2. Append the new "Now UI Dashboard" row --------------------------
$ws.Range("A11").Value = "Now UI Dashboard"
$ws.Range("B11").Value = "WebService/assets"
$ws.Range("C11").Value = "MIT"
$ws.Range("D11").Value = "https://www.creative-tim.com/product/now-ui-dashboard"
$ws.Range("E11").Value = "CSS template for web service"

# Carry over the same "last row" box formatting used by the previous
# final row (row 10) so the new row matches the rest of the table.
$ws.Range("A10:E10").Copy() | Out-Null
$ws.Range("A11:E11").PasteSpecial(-4122) | Out-Null

# --- 3. Wire up the hyperlinks -----------------------------------------
$ws.Hyperlinks.Add($ws.Range("D11"), "https://www.creative-tim.com/product/now-ui-dashboard") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://clearcanvas.github.io/") | Out-Null

# Leave the selection where the author ended up after editing C2.
$ws.Range("C2").Select() | Out-Null
